$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("StatQuery") shifting old B->C and C->D
$ws.Columns("B").Insert()

# Match new column B width to column A's (both ~75.8 "characters" wide)
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# Header for the new column
$ws.Range("B1").Value = "StatQuery"

# New stat-bar query text, matching the style (wrap text) used by A2
$statQuery = @'
MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Chesapeake Bay Retriever']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@
$ws.Range("B2").Value = $statQuery
$ws.Range("B2").WrapText = $true

# Reflect the selection state saved with the workbook: the new column B
# (whole column) is selected
$ws.Columns("B").Select()
